$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "Price" (column D) and "Volume(1h)" (column E) figures for each
# crypto row. Column D values are entered with a leading apostrophe so Excel
# keeps them as literal text (matching the source data) instead of reinterpreting
# them as numbers and silently dropping significant trailing/leading zeros.

$ws.Range("D2").Value = "'30.573.59"
$ws.Range("E2").Value = "  -0.42%  "

$ws.Range("D3").Value = "'1.883.72"
$ws.Range("E3").Value = "  -0.29%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'246.43"
$ws.Range("E5").Value = "  -0.61%  "

$ws.Range("E6").Value = "  -0.02%  "

$ws.Range("D7").Value = "'0.4730"
$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").Value = "'0.2892"
$ws.Range("E8").Value = "  -1.12%  "

$ws.Range("D10").Value = "'22.29"
$ws.Range("E10").Value = "  +1.01%  "

$ws.Range("D11").Value = "'0.7721"
$ws.Range("E11").Value = "  +4.63%  "

$ws.Range("D12").Value = "'100.81"
$ws.Range("E12").Value = "  +4.20%  "

$ws.Range("D13").Value = "'0.07800"
$ws.Range("E13").Value = "  -0.03%  "

$ws.Range("D14").Value = "'1.883.45"
$ws.Range("E14").Value = "  -0.26%  "

$ws.Range("D15").Value = "'5.257"
$ws.Range("E15").Value = "  +0.11%  "

$ws.Range("D16").Value = "'284.98"
$ws.Range("E16").Value = "  +0.18%  "

$ws.Range("D17").Value = "'30.550.52"
$ws.Range("E17").Value = "  -0.41%  "

$ws.Range("E18").Value = "  -0.33%  "

$ws.Range("D19").Value = "'0.000007522"
$ws.Range("E19").Value = "  -0.21%  "

$ws.Range("D21").Value = "'2.127.81"
$ws.Range("E21").Value = "  -0.48%  "

$ws.Range("D22").Value = "'5.354"
$ws.Range("E22").Value = "  +0.61%  "

$ws.Range("E23").Value = "  +0.21%  "

$ws.Range("D24").Value = "'6.420"
$ws.Range("E24").Value = "  +2.71%  "

$ws.Range("D25").Value = "'9.159"
$ws.Range("E25").Value = "  -0.75%  "

$ws.Range("D26").Value = "'162.41"
$ws.Range("E26").Value = "  -1.59%  "

$ws.Range("D27").Value = "'19.08"
$ws.Range("E27").Value = "  +0.65%  "

$ws.Range("D28").Value = "'1.913"
$ws.Range("E28").Value = "  -0.30%  "

$ws.Range("E29").Value = "  -0.40%  "

$ws.Range("D30").Value = "'0.09707"
$ws.Range("E30").Value = "  -0.39%  "

$ws.Range("E31").Value = "  +1.01%  "

$ws.Range("D32").Value = "'4.262"
$ws.Range("E32").Value = "  -0.89%  "

$ws.Range("D33").Value = "'4.189"
$ws.Range("E33").Value = "  -0.21%  "

$ws.Range("D34").Value = "'0.04841"
$ws.Range("E34").Value = "  -0.48%  "

$ws.Range("E35").Value = "  +0.35%  "

$ws.Range("D36").Value = "'0.6971"
$ws.Range("E36").Value = "  +0.05%  "

$ws.Range("D37").Value = "'2.760"
$ws.Range("E37").Value = "  +1.20%  "

$ws.Range("D38").Value = "'0.01914"
$ws.Range("E38").Value = "  +0.85%  "

$ws.Range("D39").Value = "'2.893"
$ws.Range("E39").Value = "  +3.09%  "

$ws.Range("D40").Value = "'76.01"
$ws.Range("E40").Value = "  -0.18%  "

$ws.Range("E41").Value = "  -0.80%  "

$ws.Range("D42").Value = "'1.982"
$ws.Range("E42").Value = "  -0.92%  "

$ws.Range("D43").Value = "'0.4251"
$ws.Range("E43").Value = "  -0.60%  "

$ws.Range("D44").Value = "'1.000"
$ws.Range("E44").Value = "  -0.06%  "

$ws.Range("D45").Value = "'0.8315"
$ws.Range("E45").Value = "  -0.51%  "

$ws.Range("D46").Value = "'101.58"
$ws.Range("E46").Value = "  -0.04%  "

$ws.Range("D47").Value = "'9.856"
$ws.Range("E47").Value = "  +3.93%  "

$ws.Range("D48").Value = "'7.027"
$ws.Range("E48").Value = "  -0.18%  "

$ws.Range("D49").Value = "'35.23"
$ws.Range("E49").Value = "  -1.02%  "

$ws.Range("D50").Value = "'892.66"
$ws.Range("E50").Value = "  -2.54%  "

$ws.Range("D51").Value = "'0.05776"
$ws.Range("E51").Value = "  +0.36%  "
